$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-01 Monday" "2025-12-02 Tuesday"

Replace-Text "762÷7=108, 6" "201÷6=33, 3"
Replace-Text "543÷4=135, 3" "218÷6=36, 2"
Replace-Text "535÷7=76, 3" "876÷9=97, 3"
Replace-Text "368÷9=40, 8" "529÷9=58, 7"
Replace-Text "253÷5=50, 3" "797÷5=159, 2"
Replace-Text "324÷2=162, 0" "721÷2=360, 1"
Replace-Text "252÷4=63, 0" "586÷2=293, 0"
Replace-Text "229÷3=76, 1" "422÷6=70, 2"
Replace-Text "222÷8=27, 6" "223÷2=111, 1"
Replace-Text "417÷9=46, 3" "618÷5=123, 3"
Replace-Text "565÷2=282, 1" "838÷4=209, 2"
Replace-Text "176÷2=88, 0" "728÷7=104, 0"
Replace-Text "415÷8=51, 7" "902÷2=451, 0"
Replace-Text "295÷6=49, 1" "135÷6=22, 3"
Replace-Text "908÷5=181, 3" "689÷8=86, 1"
Replace-Text "742÷3=247, 1" "558÷8=69, 6"
Replace-Text "986÷7=140, 6" "820÷3=273, 1"
Replace-Text "435÷4=108, 3" "147÷7=21, 0"
Replace-Text "881÷7=125, 6" "811÷8=101, 3"
Replace-Text "966÷6=161, 0" "517÷8=64, 5"
Replace-Text "631÷3=210, 1" "234÷8=29, 2"
Replace-Text "694÷3=231, 1" "935÷8=116, 7"
Replace-Text "606÷3=202, 0" "573÷9=63, 6"
Replace-Text "504÷9=56, 0" "492÷6=82, 0"
Replace-Text "450÷7=64, 2" "978÷4=244, 2"
